$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.547.84"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').Value = "'1.674.62"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.73%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'219.84"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.00%  '
$ws.Range('D6').Value = "'0.5289"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('D7').Value = "'1.002"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'0.2684"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.87%  '
$ws.Range('D9').Value = "'0.06388"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').Value = "'21.80"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.70%  '
$ws.Range('D11').Value = "'0.07801"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'4.494"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.670.29"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = "'0.5580"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').Value = "'0.0₅8332"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('D16').Value = "'65.71"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = "'26.565.59"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = "'4.771"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').Value = "'193.48"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.63%  '
$ws.Range('D21').Value = "'10.35"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').Value = "'6.317"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('D23').Value = "'1.002"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = "'0.1274"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.64%  '
$ws.Range('D25').Value = "'138.50"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.28%  '
$ws.Range('D26').Value = "'7.405"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = "'16.34"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.31%  '
$ws.Range('D28').Value = "'1.428"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +3.17%  '
$ws.Range('D29').Value = "'0.06279"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +5.68%  '
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('D31').Value = "'3.605"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +6.02%  '
$ws.Range('D32').Value = "'3.419"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').Value = "'1.692"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.75%  '
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('D35').Value = "'0.6197"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +10.43%  '
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('D38').Value = "'0.01619"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = "'6.083"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.00%  '
$ws.Range('D40').Value = "'1.095.91"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +6.37%  '
$ws.Range('D41').Value = "'0.8617"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('D42').Value = "'1.001"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').Value = "'100.61"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.35%  '
$ws.Range('D44').Value = "'1.822.52"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('E45').Value = '  +5.46%  '
$ws.Range('D46').Value = "'0.0₈109"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.17%  '
$ws.Range('D47').Value = "'8.211"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('D48').Value = "'1.534"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +11.08%  '
$ws.Range('D49').Value = "'0.9992"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').Value = "'0.05195"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('D51').Value = "'6.016"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.34%  '
